# Weekly fruit/vegetable price update ("Fruta / hortaliza, semanal").
# A new week's record is inserted at the top of the data block (row 7),
# pushing the existing data rows (old 7-39) down by one (new 8-40).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 7, shifting rows 7..39 down to 8..40.
$ws.Rows("7:7").Insert()

# Populate the newly inserted row 7 with the latest week's data.
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C7").Value = "Los Lagos"
$ws.Range("D7").Value = 44530
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 300000000
$ws.Range("G7").Value = "Espárragos"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 600
$ws.Range("K7").Value = 1500
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = 1500
$ws.Range("N7").Value = "$/kilo"
$ws.Range("O7").Value = "Provincia de Linares"
$ws.Range("P7").Value = 1500
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = "Hortaliza"
